$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 3 trailing stale rows (338-340); remaining rows shift up automatically
$ws.Rows("338:340").Delete()

# Row 149
$ws.Cells.Item(149, 2).Value = 5456603
$ws.Cells.Item(149, 6).Value = 'Lech Poznan'
$ws.Cells.Item(149, 7).Value = 'Jagiellonia Bialystok'
$ws.Cells.Item(149, 8).Value = 2
$ws.Cells.Item(149, 9).Value = 0
$ws.Cells.Item(149, 10).Value = 'H'
$ws.Cells.Item(149, 11).Value = 1.363
$ws.Cells.Item(149, 12).Value = 4.75
$ws.Cells.Item(149, 13).Value = 6.5
$ws.Cells.Item(149, 14).Value = 1.222
$ws.Cells.Item(149, 15).Value = 5.5
$ws.Cells.Item(149, 16).Value = 8
$ws.Cells.Item(149, 17).Value = -1.75
$ws.Cells.Item(149, 18).Value = 1.925
$ws.Cells.Item(149, 19).Value = 1.925
$ws.Cells.Item(149, 20).Value = 3.25
$ws.Cells.Item(149, 21).Value = 1.95
$ws.Cells.Item(149, 22).Value = 1.9
$ws.Cells.Item(149, 23).Value = 0.222
$ws.Cells.Item(149, 24).Value = -1
$ws.Cells.Item(149, 26).Value = 0.4625
$ws.Cells.Item(149, 27).Value = -0.5
$ws.Cells.Item(149, 29).Value = 0.8999999999999999

# Row 150
$ws.Cells.Item(150, 2).Value = 5428774
$ws.Cells.Item(150, 6).Value = 'Pogon Szczecin'
$ws.Cells.Item(150, 7).Value = 'Radomiak Radom'
$ws.Cells.Item(150, 8).Value = 4
$ws.Cells.Item(150, 11).Value = 1.571
$ws.Cells.Item(150, 12).Value = 4
$ws.Cells.Item(150, 13).Value = 4.75
$ws.Cells.Item(150, 14).Value = 1.533
$ws.Cells.Item(150, 15).Value = 4.333
$ws.Cells.Item(150, 16).Value = 4.75
$ws.Cells.Item(150, 17).Value = -1
$ws.Cells.Item(150, 18).Value = 1.875
$ws.Cells.Item(150, 19).Value = 1.975
$ws.Cells.Item(150, 20).Value = 3
$ws.Cells.Item(150, 21).Value = 1.875
$ws.Cells.Item(150, 22).Value = 1.975
$ws.Cells.Item(150, 23).Value = 0.5329999999999999
$ws.Cells.Item(150, 26).Value = 0.875
$ws.Cells.Item(150, 27).Value = -1
$ws.Cells.Item(150, 28).Value = 0.875
$ws.Cells.Item(150, 29).Value = -1

# Row 153
$ws.Cells.Item(153, 2).Value = 5456594
$ws.Cells.Item(153, 6).Value = 'Rakow Czestochowa'
$ws.Cells.Item(153, 7).Value = 'Zaglebie Lubin'
$ws.Cells.Item(153, 8).Value = 1
$ws.Cells.Item(153, 9).Value = 1
$ws.Cells.Item(153, 10).Value = 'D'
$ws.Cells.Item(153, 11).Value = 1.444
$ws.Cells.Item(153, 12).Value = 4.5
$ws.Cells.Item(153, 13).Value = 5.75
$ws.Cells.Item(153, 14).Value = 1.3
$ws.Cells.Item(153, 15).Value = 5.25
$ws.Cells.Item(153, 16).Value = 7
$ws.Cells.Item(153, 17).Value = -1.5
$ws.Cells.Item(153, 18).Value = 1.9
$ws.Cells.Item(153, 19).Value = 1.95
$ws.Cells.Item(153, 21).Value = 1.9
$ws.Cells.Item(153, 22).Value = 1.95
$ws.Cells.Item(153, 23).Value = -1
$ws.Cells.Item(153, 24).Value = 4.25
$ws.Cells.Item(153, 26).Value = -1
$ws.Cells.Item(153, 27).Value = 0.95
$ws.Cells.Item(153, 28).Value = -1
$ws.Cells.Item(153, 29).Value = 0.95

# Row 334
$ws.Cells.Item(334, 2).Value = 6775543
$ws.Cells.Item(334, 5).Value = 45340.35416666666
$ws.Cells.Item(334, 6).Value = 'Zaglebie Lubin'
$ws.Cells.Item(334, 7).Value = 'Cracovia Krakow'
$ws.Cells.Item(334, 11).Value = 2.4
$ws.Cells.Item(334, 12).Value = 3.4
$ws.Cells.Item(334, 13).Value = 2.8
$ws.Cells.Item(334, 14).Value = 2.375
$ws.Cells.Item(334, 15).Value = 3.4
$ws.Cells.Item(334, 16).Value = 2.9
$ws.Cells.Item(334, 18).Value = 2.025
$ws.Cells.Item(334, 19).Value = 1.825
$ws.Cells.Item(334, 20).Value = 2.5
$ws.Cells.Item(334, 21).Value = 2.025
$ws.Cells.Item(334, 22).Value = 1.825

# Row 335
$ws.Cells.Item(335, 2).Value = 6774882
$ws.Cells.Item(335, 5).Value = 45340.45833333334
$ws.Cells.Item(335, 6).Value = 'Legia Warsaw'
$ws.Cells.Item(335, 7).Value = 'MKS Puszcza Niepolomice'
$ws.Cells.Item(335, 11).Value = 1.333
$ws.Cells.Item(335, 12).Value = 5.5
$ws.Cells.Item(335, 13).Value = 7.5
$ws.Cells.Item(335, 14).Value = 1.333
$ws.Cells.Item(335, 15).Value = 5.5
$ws.Cells.Item(335, 16).Value = 7.5
$ws.Cells.Item(335, 17).Value = -1.5
$ws.Cells.Item(335, 18).Value = 2.025
$ws.Cells.Item(335, 19).Value = 1.825
$ws.Cells.Item(335, 20).Value = 2.75
$ws.Cells.Item(335, 21).Value = 1.9
$ws.Cells.Item(335, 22).Value = 1.95

# Row 336
$ws.Cells.Item(336, 2).Value = 6774460
$ws.Cells.Item(336, 5).Value = 45340.5625
$ws.Cells.Item(336, 6).Value = 'LKS Lodz'
$ws.Cells.Item(336, 7).Value = 'Widzew Lodz'
$ws.Cells.Item(336, 11).Value = 2.75
$ws.Cells.Item(336, 12).Value = 3.1
$ws.Cells.Item(336, 13).Value = 2.625
$ws.Cells.Item(336, 14).Value = 3.2
$ws.Cells.Item(336, 15).Value = 3.2
$ws.Cells.Item(336, 16).Value = 2.3
$ws.Cells.Item(336, 17).Value = 0.25
$ws.Cells.Item(336, 18).Value = 1.9
$ws.Cells.Item(336, 19).Value = 1.95
$ws.Cells.Item(336, 20).Value = 2.5
$ws.Cells.Item(336, 21).Value = 2.05
$ws.Cells.Item(336, 22).Value = 1.8

# Row 337
$ws.Cells.Item(337, 2).Value = 6775541
$ws.Cells.Item(337, 5).Value = 45341.625
$ws.Cells.Item(337, 6).Value = 'Gornik Zabrze'
$ws.Cells.Item(337, 7).Value = 'Korona Kielce'
$ws.Cells.Item(337, 11).Value = 2.15
$ws.Cells.Item(337, 12).Value = 3.2
$ws.Cells.Item(337, 13).Value = 3.6
$ws.Cells.Item(337, 14).Value = 2.05
$ws.Cells.Item(337, 15).Value = 3.25
$ws.Cells.Item(337, 16).Value = 3.75
$ws.Cells.Item(337, 17).Value = -0.5
$ws.Cells.Item(337, 18).Value = 2.1
$ws.Cells.Item(337, 19).Value = 1.775
$ws.Cells.Item(337, 21).Value = 2.05
$ws.Cells.Item(337, 22).Value = 1.8

